$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set the Runmode column (C) to "Y" for all test suite rows (running all the suites)
$ws.Range("C3:C7").Value = "Y"

# Update the selection to reflect the reviewed range
$ws.Range("C2:C7").Select()
